# "some error fixes to penAndPaper 2"
# The edit corrects the a(X)/b(X) sample values on worksheet "2" (row 6),
# which feed the S(X) similarity formulas in columns L:R, and leaves
# worksheet "2" as the active/selected sheet (with R6 selected).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2")

# Fix the a(X) values (B6:E6) - B6 stays, C6/D6/E6 corrected
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 4.49535804129924
$ws.Range("E6").Value = 4.9286503810670403

# Fix the b(X) values (G6:J6) - G6 stays, H6/I6/J6 corrected
$ws.Range("H6").Value = 6.9823759942501198
$ws.Range("I6").Value = 6
$ws.Range("J6").Value = 6.4031242374328396

# Make worksheet "2" the active sheet/tab, with R6 selected
$ws.Activate()
$ws.Range("R6").Select()

# Recalculate so the S(X) formulas (L6:R6) pick up the corrected inputs
$excel.Calculate()
